# daily auto push: 2026-02-16 22:47 UTC
# Insert the next daily data point at row 808 (2026/02/17, weekday 火,
# hour 5, ranking 57), pushing the existing 2026/12/29..2027/01/05 rows
# down by one (808->809 ... 849->850).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 808:849 down to 809:850, leaving a blank row 808.
$ws.Rows.Item(808).Insert()

# The "日付" column stores dates as plain text (e.g. "2026/12/29"), not
# real Excel dates, so force text format before assigning the value to
# stop Excel from auto-converting the string to a date serial number.
$ws.Range("A808").NumberFormat = "@"
$ws.Range("A808").Value = "2026/02/17"
# Drop the temporary text-format override so the cell keeps the same
# (default) style as every other data row.
$ws.Range("A808").ClearFormats()

$ws.Range("B808").Value = "火"
$ws.Range("C808").Value = 5
$ws.Range("D808").Value = 57
